# Auto-generated script applying scheduled market-data refresh to Sheets/Siren_Profits.xlsx
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per leve row across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 965
$ws.Range("I20").Value = 965
$ws.Range("K20").Value = 965
$ws.Range("M20").Value = -735

$ws.Range("H35").Value = 965
$ws.Range("I35").Value = 965
$ws.Range("K35").Value = 965
$ws.Range("M35").Value = -586

$ws.Range("H39").Value = 464.64285
$ws.Range("J39").Value = 741
$ws.Range("L39").Value = 2223
$ws.Range("N39").Value = -2815

$ws.Range("H80").Value = 161501.77
$ws.Range("I80").Value = 253825.5
$ws.Range("J80").Value = 13783.8
$ws.Range("K80").Value = 761476.5
$ws.Range("L80").Value = 41351.39999999999
$ws.Range("M80").Value = -760478.5
$ws.Range("N80").Value = -43347.39999999999

$ws.Range("H83").Value = 161501.77
$ws.Range("I83").Value = 253825.5
$ws.Range("J83").Value = 13783.8
$ws.Range("K83").Value = 2284429.5
$ws.Range("L83").Value = 124054.2
$ws.Range("M83").Value = -2279437.5
$ws.Range("N83").Value = -134038.2

$ws.Range("H107").Value = 5785.5356
$ws.Range("I107").Value = 6167.9473
$ws.Range("K107").Value = 6167.9473
$ws.Range("M107").Value = -4247.9473

$ws.Range("H112").Value = 57141.1
$ws.Range("J112").Value = 81069.42999999999
$ws.Range("L112").Value = 243208.29
$ws.Range("N112").Value = -245424.29

$ws.Range("H113").Value = 25782.166
$ws.Range("I113").Value = 25782.166
$ws.Range("K113").Value = 25782.166
$ws.Range("M113").Value = -22528.166

$ws.Range("H116").Value = 658306.9399999999
$ws.Range("I116").Value = 1591139.6
$ws.Range("K116").Value = 1591139.6
$ws.Range("M116").Value = -1587697.6

$ws.Range("H131").Value = 6770.697
$ws.Range("I131").Value = 1781.7
$ws.Range("K131").Value = 5345.1
$ws.Range("M131").Value = -305.1000000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1824.5
$ws.Range("I110").Value = 1777
$ws.Range("J110").Value = 1998.6666
$ws.Range("K110").Value = 1777
$ws.Range("L110").Value = 1998.6666
$ws.Range("M110").Value = 268
$ws.Range("N110").Value = -6088.6666

$ws.Range("H122").Value = 1278268
$ws.Range("I122").Value = 4207
$ws.Range("J122").Value = 3507874.8
$ws.Range("K122").Value = 12621
$ws.Range("L122").Value = 10523624.4
$ws.Range("M122").Value = -10171
$ws.Range("N122").Value = -10528524.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1015.5714
$ws.Range("I80").Value = 1084.6
$ws.Range("J80").Value = 952.8182
$ws.Range("K80").Value = 1084.6
$ws.Range("L80").Value = 952.8182
$ws.Range("M80").Value = -86.59999999999991
$ws.Range("N80").Value = -2948.8182

$ws.Range("H83").Value = 1015.5714
$ws.Range("I83").Value = 1084.6
$ws.Range("J83").Value = 952.8182
$ws.Range("K83").Value = 5423
$ws.Range("L83").Value = 4764.091
$ws.Range("M83").Value = -431
$ws.Range("N83").Value = -14748.091

$ws.Range("H134").Value = 12003.866
$ws.Range("I134").Value = 21401.715
$ws.Range("K134").Value = 64205.145
$ws.Range("M134").Value = -61670.145

$ws.Range("H140").Value = 87949
$ws.Range("J140").Value = 87949
$ws.Range("L140").Value = 87949
$ws.Range("N140").Value = -98309

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4019.5
$ws.Range("I16").Value = 4019.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4019.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3732.5
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 3788.389
$ws.Range("I31").Value = 2747.1
$ws.Range("K31").Value = 2747.1
$ws.Range("M31").Value = -2452.1

$ws.Range("H34").Value = 3788.389
$ws.Range("I34").Value = 2747.1
$ws.Range("K34").Value = 2747.1
$ws.Range("M34").Value = -2545.1

$ws.Range("H99").Value = 14521888
$ws.Range("I99").Value = 14521888
$ws.Range("K99").Value = 14521888
$ws.Range("M99").Value = -14520390

$ws.Range("H113").Value = 4019.5
$ws.Range("I113").Value = 4019.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4019.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1849.5
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 2306.375
$ws.Range("I122").Value = 2188
$ws.Range("K122").Value = 6564
$ws.Range("M122").Value = -4114

$ws.Range("H126").Value = 14521888
$ws.Range("I126").Value = 14521888
$ws.Range("K126").Value = 43565664
$ws.Range("M126").Value = -43563194

$ws.Range("H134").Value = 2724543.8
$ws.Range("I134").Value = 2983609.2
$ws.Range("K134").Value = 8950827.600000001
$ws.Range("M134").Value = -8948292.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 6796.3
$ws.Range("I112").Value = 6192.6
$ws.Range("J112").Value = 7400
$ws.Range("K112").Value = 18577.8
$ws.Range("L112").Value = 22200
$ws.Range("M112").Value = -17469.8
$ws.Range("N112").Value = -24416

$ws.Range("H129").Value = 25642538
$ws.Range("J129").Value = 55558220
$ws.Range("L129").Value = 166674660
$ws.Range("N129").Value = -166684660

$ws.Range("H140").Value = 12194.8125
$ws.Range("I140").Value = 12194.8125
$ws.Range("K140").Value = 36584.4375
$ws.Range("M140").Value = -31404.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 942390.8
$ws.Range("J49").Value = 942390.8
$ws.Range("L49").Value = 942390.8
$ws.Range("N49").Value = -942758.8

$ws.Range("H59").Value = 14750
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 19500
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 19500
$ws.Range("M59").Value = -9417
$ws.Range("N59").Value = -20666

$ws.Range("H80").Value = 3083.7778
$ws.Range("I80").Value = 2989.5
$ws.Range("J80").Value = 3159.2
$ws.Range("K80").Value = 2989.5
$ws.Range("L80").Value = 3159.2
$ws.Range("M80").Value = -1991.5
$ws.Range("N80").Value = -5155.2

$ws.Range("H83").Value = 3083.7778
$ws.Range("I83").Value = 2989.5
$ws.Range("J83").Value = 3159.2
$ws.Range("K83").Value = 14947.5
$ws.Range("L83").Value = 15796
$ws.Range("M83").Value = -9955.5
$ws.Range("N83").Value = -25780

$ws.Range("H101").Value = 32277
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").ClearContents()

$ws.Range("H102").Value = 10801.5
$ws.Range("I102").Value = 11729.454
$ws.Range("J102").Value = 594
$ws.Range("K102").Value = 11729.454
$ws.Range("L102").Value = 594
$ws.Range("M102").Value = -10107.454
$ws.Range("N102").Value = -3838

$ws.Range("H113").Value = 2430.5386
$ws.Range("I113").Value = 2633.3333
$ws.Range("J113").Value = 2256.7144
$ws.Range("K113").Value = 2633.3333
$ws.Range("L113").Value = 2256.7144
$ws.Range("M113").Value = -463.3332999999998
$ws.Range("N113").Value = -6596.7144

$ws.Range("H122").Value = 13839
$ws.Range("I122").Value = 16118.117
$ws.Range("K122").Value = 48354.351
$ws.Range("M122").Value = -45904.351

$ws.Range("H123").Value = 19230.691
$ws.Range("J123").Value = 19230.691
$ws.Range("L123").Value = 19230.691
$ws.Range("N123").Value = -24130.691

$ws.Range("H126").Value = 16594.76
$ws.Range("I126").Value = 16626.084
$ws.Range("K126").Value = 49878.25199999999
$ws.Range("M126").Value = -47408.25199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 71079.17999999999
$ws.Range("I40").Value = 81986.78
$ws.Range("K40").Value = 81986.78
$ws.Range("M40").Value = -81850.78

$ws.Range("H46").Value = 2345.4285
$ws.Range("J46").Value = 2754.05
$ws.Range("L46").Value = 2754.05
$ws.Range("N46").Value = -3130.05

$ws.Range("H93").Value = 5395.5
$ws.Range("I93").Value = 5921.8184
$ws.Range("K93").Value = 5921.8184
$ws.Range("M93").Value = -4673.8184

$ws.Range("H122").Value = 4638.9
$ws.Range("I122").Value = 4440.25
$ws.Range("K122").Value = 13320.75
$ws.Range("M122").Value = -10870.75

$ws.Range("H132").Value = 1157249.2
$ws.Range("I132").Value = 2144376.5
$ws.Range("K132").Value = 6433129.5
$ws.Range("M132").Value = -6430599.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 5000
$ws.Range("I96").Value = 5000
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 5000
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -3627
$ws.Range("N96").ClearContents()

$ws.Range("H122").Value = 2412.4814
$ws.Range("I122").Value = 1832.3043
$ws.Range("K122").Value = 5496.9129
$ws.Range("M122").Value = -3046.9129

$ws.Range("H126").Value = 22764.953
$ws.Range("I126").Value = 25769.059
$ws.Range("K126").Value = 77307.177
$ws.Range("M126").Value = -74837.177

$ws.Range("H132").Value = 19469.637
$ws.Range("J132").Value = 10297.833
$ws.Range("L132").Value = 30893.499
$ws.Range("N132").Value = -35953.499
